# All object types - dynamic identification of keys
# Add a new worksheet "AllObjectTypes" after the last existing sheet,
# mirroring the layout of "CreateComplaint" but pointing at the new
# "validateAllObjTypes" test data.

$wb = $excel.ActiveWorkbook

$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws = $wb.Worksheets.Add($null, $lastSheet)
$ws.Name = "AllObjectTypes"

$ws.Range("A1").Value = "validateAllObjTypes"
$ws.Range("B1").Value = "key"
$ws.Range("B2").Value = "value"

$ws.Range("D13").Select()
